$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Candidate ID 23091337 -> 23091510)
$ws.Range("A2").Value = "vCmwh683"
$ws.Range("B2").Value = 23091510
$ws.Range("C2").Value = "itdccrx94"
$ws.Range("D2").Value = "BaY3`$6&t"
$ws.Range("F2").Value = "RnxoWuxJ"
$ws.Range("G2").Value = "dSKD"

# Row 3 (Candidate ID 23091336 -> 23091509)
$ws.Range("A3").Value = "zkFjv156"
$ws.Range("B3").Value = 23091509
$ws.Range("C3").Value = "rhqqhmo10"
$ws.Range("D3").Value = "G&Ps4!u9"
$ws.Range("F3").Value = "ArKJlAuf"
$ws.Range("G3").Value = "ymuy"

# Row 4 (Candidate ID 23091335 -> 23091508)
$ws.Range("A4").Value = "yHmgA454"
$ws.Range("B4").Value = 23091508
$ws.Range("C4").Value = "vpyvkcc63"
$ws.Range("D4").Value = "qr&4`$WM3"
$ws.Range("F4").Value = "xrUDENCY"
$ws.Range("G4").Value = "iLMK"
